# Applies the "disabling fields according to other field inputs" edit
# to the Piql partner order form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields -----------------------------------------------------
# Date/time the order form was filled in (same day, later time).
$ws.Range("G4").Value = 44080.8674935521

# Customer name changed.
$ws.Range("G7").Value = "El caballo vengador"

# Address "comment" field changed.
$ws.Range("F10").Value = "Tirando las maletas"

# --- Line items ----------------------------------------------------------

# Row 18: piqlConnect (only piqlFilm) - now filled in (Qty/Unit price/Price)
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1500
$ws.Range("H18").Value = 1500

# Row 19: Digital (GB) - quantity and total price updated
$ws.Range("F19").Value = 450
$ws.Range("H19").Value = 6750

# Row 21: Online Storage section header - cleared out
$ws.Range("F21").Value = $null
$ws.Range("G21").Value = $null
$ws.Range("H21").Value = $null

# Row 22: Online Storage (GB) - payment type + values cleared out
$ws.Range("E22").Value = $null
$ws.Range("F22").Value = $null
$ws.Range("G22").Value = $null
$ws.Range("H22").Value = $null

# Row 24: Registration fee - cleared out
$ws.Range("F24").Value = $null
$ws.Range("G24").Value = $null
$ws.Range("H24").Value = $null

# Row 25: AWA contribution - entity + values cleared out
$ws.Range("E25").Value = $null
$ws.Range("F25").Value = $null
$ws.Range("G25").Value = $null
$ws.Range("H25").Value = $null

# Row 26: Management fee (per period) - cleared out
$ws.Range("F26").Value = $null
$ws.Range("G26").Value = $null
$ws.Range("H26").Value = $null

# Row 27: Storage (per reel / per period) - period + values cleared out
$ws.Range("E27").Value = $null
$ws.Range("F27").Value = $null
$ws.Range("G27").Value = $null
$ws.Range("H27").Value = $null

# Row 28: Professional Services (per day) - cleared out
$ws.Range("F28").Value = $null
$ws.Range("G28").Value = $null
$ws.Range("H28").Value = $null

# Row 29: piqlReader - cleared out
$ws.Range("F29").Value = $null
$ws.Range("G29").Value = $null
$ws.Range("H29").Value = $null

# Row 30: Installation and training - cleared out
$ws.Range("F30").Value = $null
$ws.Range("G30").Value = $null
$ws.Range("H30").Value = $null

# Row 31: Service agreement (per year) - type + values cleared out
$ws.Range("E31").Value = $null
$ws.Range("F31").Value = $null
$ws.Range("G31").Value = $null
$ws.Range("H31").Value = $null

# Row 32: Shipment cost - reels count and resulting price changed
$ws.Range("E32").Value = 4
$ws.Range("G32").Value = 30
$ws.Range("H32").Value = 120

# --- Totals ---------------------------------------------------------------
$ws.Range("H33").Value = 8250
$ws.Range("H34").Value = 0
